$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 — empty cell, same style as Q2 (thin bottom-bordered blank)
$ws.Range("Q2").Copy()
$ws.Range("R2").PasteSpecial(-4122)

# Row 3 — year header 2021, same style as Q3
$ws.Range("Q3").Copy()
$ws.Range("R3").PasteSpecial(-4122)
$ws.Range("R3").Value = 2021

# Row 4 — new data point, same style as P4/Q4 plus a "0.0" number format
$ws.Range("P4").Copy()
$ws.Range("R4").PasteSpecial(-4122)
$ws.Range("R4").Value = 18
$ws.Range("R4").NumberFormat = "0.0"

# Row 5 — same style as Q5
$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)
$ws.Range("R5").Value = 1.7480265877296817

# Rows 6-12 — copy style from column D (style 25), matching how these rows were
# originally extended (column Q on these rows carries a different, bottom-less
# style than the rest of the row)
$ws.Range("D6").Copy()
$ws.Range("R6").PasteSpecial(-4122)
$ws.Range("R6").Value = 4.1112601249414027

$ws.Range("D7").Copy()
$ws.Range("R7").PasteSpecial(-4122)
$ws.Range("R7").Value = 1.5225742120245318

$ws.Range("D8").Copy()
$ws.Range("R8").PasteSpecial(-4122)
$ws.Range("R8").Value = 1.2326518235454269

$ws.Range("D9").Copy()
$ws.Range("R9").PasteSpecial(-4122)
$ws.Range("R9").Value = 4.0865392096984241

$ws.Range("D10").Copy()
$ws.Range("R10").PasteSpecial(-4122)
$ws.Range("R10").Value = 0.84876624403485645

$ws.Range("D11").Copy()
$ws.Range("R11").PasteSpecial(-4122)
$ws.Range("R11").Value = 2.1456657699653627

$ws.Range("D12").Copy()
$ws.Range("R12").PasteSpecial(-4122)
$ws.Range("R12").Value = 1.8214779402142154

# Row 13 — bottom totals row, same style as Q13
$ws.Range("Q13").Copy()
$ws.Range("R13").PasteSpecial(-4122)
$ws.Range("R13").Value = 0.51989507542472779

# Update the active sheet view's selection to match the post-edit state
$ws.Range("R24:R25").Select()
